$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set A2 to the "MONTEREY AREA TOTALS" label (moved from B2)
$ws.Range("A2").Value = "MONTEREY AREA TOTALS"

# Set B2 to "Totals"
$ws.Range("B2").Value = "Totals"

# Update the selection to B3
$ws.Range("B3").Select()
